# Bot5 GUI testeada 260123
# Flip the "parametrosInicio" start-parameter cell (B12) from 1 to 0,
# and leave the sheet's selection sitting on that same cell, matching
# the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parametrosInicio")
$ws.Activate()

$cell = $ws.Range("B12")
$cell.Value = 0
$cell.Select() | Out-Null
